$d = $word.ActiveDocument

# Locate the paragraph that ends the document body ("Je cherche comment ...
# faisable déjà !") so the new bullet is appended right after it.
$anchorText = "faisable d"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$anchorText*") {
        $target = $p
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Item($d.Paragraphs.Count)
}

$rng = $target.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# The newly created paragraph inherits the list style/numbering
# (Paragraphedeliste, numId 1) from the paragraph it follows.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "Apparemment c'est chaud patate…"
